$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 415.5232050914462

# Row 3
$ws.Range("B3").Value = 0.006876353814593728
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 249.35809635783

# Row 4
$ws.Range("B4").Value = 0.06328177979961902
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 116886.6739907443
$ws.Range("E4").Value = 71517.89157740913
$ws.Range("G4").Value = 188406.2820863922

# Row 5
$ws.Range("B5").Value = 0.006876353814593728
$ws.Range("C5").Value = 0.3375848360084654
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 9.908488693797258
